$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift each "Week N" header along row 2 forward by one week.
# The header cells are the first cell of each merged 7-day block.
$headerCells = @("B2", "I2", "P2", "W2", "AD2", "AK2", "AR2", "AY2", "BF2", "BM2")
$weekNumbers = @(15, 16, 17, 18, 19, 20, 21, 22, 23, 24)

for ($i = 0; $i -lt $headerCells.Length; $i++) {
    $ws.Range($headerCells[$i]).Value = "Week " + $weekNumbers[$i]
}

# Update the active selection to match the edited workbook.
$ws.Range("BM3").Select()
